$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 205 - this shifts the existing rows 205..241
# down to 206..242 (matching the diff, which is a weekly-data insertion
# for "Feria Lagunitas de Puerto Montt - Pepino ensalada").
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new weekly record.
$ws.Range("A205").Value = 4
$ws.Range("B205").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value = "Los Lagos"
$ws.Range("D205").Value = 44637
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = 100112043
$ws.Range("G205").Value = "Pepino ensalada"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 200
$ws.Range("K205").Value = 24000
$ws.Range("L205").Value = 24000
$ws.Range("M205").Value = 24000
$ws.Range("N205").Value = "$/caja 60 unidades"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 400
$ws.Range("Q205").Value = 60
$ws.Range("R205").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the table.
$ws.Range("D205").NumberFormat = $ws.Range("D206").NumberFormat
